# Apply the "cryptos list" refresh: Wed Feb 22 14:47:45 UTC 2023 GitHub Actions run.
# Updates Price (column D) and Volume(1h) (column E) for each coin row, and
# for rows 35/37 also swaps the ImmutableX / InternetComputer(DFINITY) entries
# (Coin name + Link) to their corrected positions.
#
# Numeric-looking Price strings (e.g. "1.002", "0.6920") are written with a
# leading apostrophe so Excel keeps them as literal text (matching the source
# data, which stores every Price/Volume cell as text) instead of silently
# re-parsing them into numbers and dropping significant trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.111.81'
$ws.Range("E2").Value = '  -2.64%  '

# Row 3
$ws.Range("D3").Value = '1.640.08'
$ws.Range("E3").Value = '  -2.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").Value = '''309.23'
$ws.Range("E5").Value = '  -1.65%  '

# Row 6
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("D7").Value = '''0.3933'
$ws.Range("E7").Value = '  +0.35%  '

# Row 8
$ws.Range("D8").Value = '''0.3862'
$ws.Range("E8").Value = '  -2.43%  '

# Row 9
$ws.Range("D9").Value = '''1.002'
$ws.Range("E9").Value = '  +0.25%  '

# Row 10
$ws.Range("D10").Value = '''50.21'
$ws.Range("E10").Value = '  -3.44%  '

# Row 11
$ws.Range("D11").Value = '''1.365'
$ws.Range("E11").Value = '  -2.78%  '

# Row 12
$ws.Range("D12").Value = '''0.08562'
$ws.Range("E12").Value = '  -1.11%  '

# Row 13
$ws.Range("D13").Value = '''23.72'
$ws.Range("E13").Value = '  -6.36%  '

# Row 14
$ws.Range("D14").Value = '''7.080'
$ws.Range("E14").Value = '  -3.43%  '

# Row 15
$ws.Range("D15").Value = '''0.00001284'
$ws.Range("E15").Value = '  -2.72%  '

# Row 16
$ws.Range("D16").Value = '''7.498'
$ws.Range("E16").Value = '  -3.61%  '

# Row 17
$ws.Range("D17").Value = '1.632.71'
$ws.Range("E17").Value = '  -2.82%  '

# Row 18
$ws.Range("D18").Value = '''93.73'
$ws.Range("E18").Value = '  +0.02%  '

# Row 19
$ws.Range("D19").Value = '''0.06914'
$ws.Range("E19").Value = '  -2.19%  '

# Row 20
$ws.Range("D20").Value = '''20.32'
$ws.Range("E20").Value = '  -0.07%  '

# Row 21
$ws.Range("D21").Value = '''6.928'
$ws.Range("E21").Value = '  -2.10%  '

# Row 22
$ws.Range("D22").Value = '''1.002'
$ws.Range("E22").Value = '  +0.05%  '

# Row 23
$ws.Range("D23").Value = '''13.61'
$ws.Range("E23").Value = '  -2.41%  '

# Row 24
$ws.Range("D24").Value = '24.127.04'
$ws.Range("E24").Value = '  -2.58%  '

# Row 25
$ws.Range("D25").Value = '''2.404'
$ws.Range("E25").Value = '  +2.42%  '

# Row 26
$ws.Range("D26").Value = '''2.881'
$ws.Range("E26").Value = '  +2.83%  '

# Row 27
$ws.Range("D27").Value = '''22.24'
$ws.Range("E27").Value = '  -5.04%  '

# Row 28
$ws.Range("D28").Value = '''158.28'
$ws.Range("E28").Value = '  -2.62%  '

# Row 29
$ws.Range("D29").Value = '''139.97'
$ws.Range("E29").Value = '  -5.16%  '

# Row 30
$ws.Range("D30").Value = '''8.093'
$ws.Range("E30").Value = '  +2.43%  '

# Row 31
$ws.Range("D31").Value = '''5.271'
$ws.Range("E31").Value = '  -9.77%  '

# Row 32
$ws.Range("D32").Value = '''2.477'
$ws.Range("E32").Value = '  +2.83%  '

# Row 33
$ws.Range("D33").Value = '1.813.03'
$ws.Range("E33").Value = '  -3.02%  '

# Row 34
$ws.Range("D34").Value = '''0.08066'
$ws.Range("E34").Value = '  -4.32%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''0.9710'
$ws.Range("E35").Value = '  -3.09%  '

# Row 36
$ws.Range("E36").Value = '  -4.99%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''6.696'
$ws.Range("E37").Value = '  -4.06%  '

# Row 38
$ws.Range("D38").Value = '''0.2693'
$ws.Range("E38").Value = '  -4.50%  '

# Row 39
$ws.Range("D39").Value = '''0.09233'
$ws.Range("E39").Value = '  -2.61%  '

# Row 40
$ws.Range("D40").Value = '''10.38'
$ws.Range("E40").Value = '  -1.88%  '

# Row 41
$ws.Range("D41").Value = '''1.429'
$ws.Range("E41").Value = '  -5.32%  '

# Row 42
$ws.Range("D42").Value = '''0.7528'
$ws.Range("E42").Value = '  -5.42%  '

# Row 43
$ws.Range("D43").Value = '''13.10'
$ws.Range("E43").Value = '  -3.66%  '

# Row 44
$ws.Range("D44").Value = '''16.19'
$ws.Range("E44").Value = '  -2.59%  '

# Row 45
$ws.Range("D45").Value = '''0.6920'

# Row 46
$ws.Range("D46").Value = '''2.459'
$ws.Range("E46").Value = '  -4.49%  '

# Row 47
$ws.Range("D47").Value = '''4.091'
$ws.Range("E47").Value = '  -2.34%  '

# Row 48
$ws.Range("D48").Value = '''1.002'
$ws.Range("E48").Value = '  -0.02%  '

# Row 49
$ws.Range("D49").Value = '''0.08351'
$ws.Range("E49").Value = '  -3.99%  '

# Row 50
$ws.Range("D50").Value = '''1.265'
$ws.Range("E50").Value = '  -6.36%  '

# Row 51
$ws.Range("D51").Value = '''133.42'
$ws.Range("E51").Value = '  -3.45%  '
